# Auto-update hourly job matches and history

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-10 (Title, Company, Location, Match Score (%), Matched Keywords, Posted At, Apply Link)
$data = @(
    @("Sr Software Engineer - Content Platform Engineering", "nan", "Glendale, CA, US USA", 18.9, "Data Scientist, RAG, Prompt Engineering, TensorFlow, PyTorch, S3, EC2, Docker, Kubernetes, Jenkins", "2026-02-24", "https://www.indeed.com/viewjob?jk=2dc4184d266d0070"),
    @("Sr Machine Learning Engineer", "nan", "Seattle, WA, US USA", 17.8, "Data Scientist, Machine Learning Engineer, RAG, LLaMA, Gemini, TensorFlow, PyTorch, AWS SageMaker, Git, Databricks", "2026-02-24", "https://www.indeed.com/viewjob?jk=03889315af72efca"),
    @("Sr Machine Learning Engineer", "The Walt Disney Company", "Lake Buena Vista, FL, US USA", 13.3, "AI Engineer, Machine Learning Engineer, RAG, TensorFlow, PyTorch, Docker, Kubernetes, Git, Python, R", "2026-02-24", "https://www.indeed.com/viewjob?jk=296b70730c29edd4"),
    @("Senior Backend Engineer", "Glassbox", "Remote, US USA", 12.2, "Data Scientist, Copilot, Kubernetes, Git, Kafka, Cassandra, NoSQL, SQL, R, Java", "2026-02-24", "https://www.indeed.com/viewjob?jk=a7eed342df94f18f"),
    @("Data Scientist", "Ascendion", "Chicago, IL, US USA", 11.1, "Data Scientist, FAISS, TensorFlow, PyTorch, XGBoost, Git, Hadoop, Python, R, Scala", "2026-02-24", "https://www.indeed.com/viewjob?jk=24526c3c91ea4416"),
    @("Specialist - Architecture", "LTIMindtree", "Irving, TX, US USA", 10, "RAG, Docker, Kubernetes, Kafka, Python, SQL, R, Java, Scala", "2026-02-24", "https://www.indeed.com/viewjob?jk=ef5a7c704ef11e31"),
    @("Senior Software Engineer", "CNH Industrial", "Oak Brook, IL, US USA", 10, "RAG, Docker, Kubernetes, CI/CD, Git, NoSQL, SQL, R, Scala", "2026-02-24", "https://www.indeed.com/viewjob?jk=0d56bb5a0c74a4a6"),
    @("Data Scientist", "Partify Inc.", "Warren, MI, US USA", 10, "Data Scientist, RAG, Power BI, Python, SQL, R, Scala, Optimization, Hypothesis Testing", "2026-02-24", "https://www.indeed.com/viewjob?jk=fac62c8a9b25af94"),
    @("Sr Software Engineer", "Disney Experiences", "Orlando, FL, US USA", 10, "RAG, Docker, Terraform, Git, NoSQL, SQL, R, Java, Scala", "2026-02-24", "https://www.indeed.com/viewjob?jk=ad9da3683bfea092")
)

# Column F holds date-look-alike text (e.g. "2026-02-24"); force text
# formatting so Excel does not auto-convert it to a date serial number.
$ws.Range("F2:F10").NumberFormat = "@"

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

# Remove rows 11-13 (no longer present after the update) so the sheet
# dimension shrinks from A1:G13 down to A1:G10
$ws.Range("A11:G13").Delete()

$wb.Save()
